# pso - mvo till 8 jan2016-7jan2021
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the most recent "Date range" block (rows 31-36) into a new
# block at rows 38-43, preserving styles / merged-cell layout / row heights.
$ws.Range("A31:AD36").Copy($ws.Range("A38"))

# The whole-block copy stamps the merged B-column anchor style (s=12/7/10)
# across the whole merge; restore the individual "bottom half" cell style
# (s=3) that the source rows 32/34/36 actually carry.
$ws.Range("B32").Copy()
$ws.Range("B39").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("B34").Copy()
$ws.Range("B41").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("B36").Copy()
$ws.Range("B43").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Rows 42/43 (mirroring 35/36) should only carry label cells, not the
# blank numeric placeholders the whole-block copy produced.
$ws.Range("D42:AD43").ClearContents()

# New date-range label for the appended block.
$ws.Range("A38").Value = "8-Jan-2016 :: 7-Jan-2021"

# New PSO row (38) values.
$row38 = @(29.3835376580014, 17.719533768604, 1.45621989804951, 2.8401587000681, 6.9628007630853, 4.28442043112754, 2.55465963998268, 7.82061704303467, 0.14935983281798, 9.04523013960718, 2.4605414793341, 2.63691518133205, 5.35161980724609, 1.11909785347891, 5.92994158749161, 2.38766071009967, 1.45926160769788, 3.74434555777137, 0.915351208217588, 9.44395587802283, 8.46058118583972, 7.66256398312449, 0.1665872738048, 0.603758506956849, 6.65280735056821, 1.92775944360194, 5.42000483568845)
for ($i = 0; $i -lt $row38.Length; $i++) {
    $ws.Cells.Item(38, 4 + $i).Value = $row38[$i]
}

# New Sortino row (39) values.
$row39 = @(30.5268165699694, 11.616707475663, 2.31966042240653, 6.02249276516474, 8.47720223096306, 4.85283424798499, 3.84056855186565, 2.49153965834089, 0.670632756246172, 1.24747858990165, 5.4704593694528, 4.32524387709977, 7.48709105822475, 0.707589925219238, 2.0487872696163, 4.04848964341325, 0.0904522709183136, 8.08381615795677, 3.96798421868443, 8.79700208858593, 4.9705888829398, 6.41866648844545, 4.60116072087233, 0.518750740424009, 7.9360813274762, 0.867491217557508, 2.05759594264602)
for ($i = 0; $i -lt $row39.Length; $i++) {
    $ws.Cells.Item(39, 4 + $i).Value = $row39[$i]
}

# New ACO/MVO row (40) values.
$row40 = @(28.8, 17.8885438199983, 1.40984085981362, 0.227509410207713, 11.7399060855447, 0.580882923558929, 8.96588452417452, 2.93697282019851, 2.09668392109204, 8.82697051487464, 18.5455319755054, 0.961675586014997, 7.71002905358488, 0.129010244836969, 0.0133902191114465, 2.93204243887695, 4.10046519946671, 7.69847090485439, 1.4748551806485, 0.348560308015893, 4.12033278602222, 1.75960121490938, 0.292881170523793, 1.06304360597416, 11.5268607735217, 1.09245161713521, 0.855987521346332)
for ($i = 0; $i -lt $row40.Length; $i++) {
    $ws.Cells.Item(40, 4 + $i).Value = $row40[$i]
}

# New ACO/Sortino row (41) values.
$row41 = @(30.2, 11.6275534829989, 2.28938959850171, 2.6569208667856, 7.45181052796176, 10.7459695601984, 2.44797199903069, 1.3023077053619, 5.11600331728191, 6.99989873158132, 0.972635935450368, 1.90716900911576, 17.3964788549832, 0.348705270320119, 1.07274093988888, 0.648931748930941, 3.47701687007473, 9.44086922459458, 0.889158989965562, 0.299414480173068, 0.226735822947277, 7.64815173670969, 1.12874807800403, 0.612597631014368, 14.2502716919244, 0.630450079369199, 2.3290409283321)
for ($i = 0; $i -lt $row41.Length; $i++) {
    $ws.Cells.Item(41, 4 + $i).Value = $row41[$i]
}

# Update the view to match the edited state (new active cell / scroll position).
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("D21").Select()
